# Elevator Pitch edit: rewrite the narrative paragraph with the
# Fall-2018-deployment framing, add a tracked-change placeholder for the
# course title (author: Hadfield, Steven M Civ USAF USAFA USAFA/DFCS),
# and drop the stray _GoBack bookmark from the trailing empty paragraph.

$d = $word.ActiveDocument
$word.UserName = "Hadfield, Steven M Civ USAF USAFA USAFA/DFCS"

$targetXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Last </w:t></w:r><w:r><w:t>f</w:t></w:r><w:r><w:t xml:space="preserve">all, the </w:t></w:r><w:r><w:t>Military and Strategic Studies D</w:t></w:r><w:r><w:t xml:space="preserve">epartment had a team of cadets developing a strategic wargame for the new MSS 251 </w:t></w:r><w:ins w:id="0" w:author="Hadfield, Steven M Civ USAF USAFA USAFA/DFCS" w:date="2018-09-10T17:20:00Z"><w:r><w:rPr><w:highlight w:val="yellow"/><w:rPrChange w:id="1" w:author="Hadfield, Steven M Civ USAF USAFA USAFA/DFCS" w:date="2018-09-10T17:30:00Z"><w:rPr/></w:rPrChange></w:rPr><w:t>&lt;ADD COURSE TITLE HERE&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:ins><w:r><w:t xml:space="preserve">course. They made a board game with 3D printed pieces and a typed-up rule book. The board game worked well for the test sections in </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Spring</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 2019</w:t></w:r><w:r><w:t xml:space="preserve">, but they needed </w:t></w:r><w:r><w:t>a more scalable</w:t></w:r><w:r><w:t xml:space="preserve"> platform </w:t></w:r><w:r><w:t xml:space="preserve">to support </w:t></w:r><w:bookmarkStart w:id="2" w:name="_GoBack"/><w:bookmarkEnd w:id="2"/><w:r><w:t xml:space="preserve">the 27 live sections for the next Fall, so they contacted the </w:t></w:r><w:r><w:t>D</w:t></w:r><w:r><w:t>epartment</w:t></w:r><w:r><w:t xml:space="preserve"> of Computer and Cyber Sciences</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">Cadet First Class </w:t></w:r><w:r><w:t xml:space="preserve">Jack Kulp </w:t></w:r><w:r><w:t>began development of</w:t></w:r><w:r><w:t xml:space="preserve"> a browser-based version of the board game </w:t></w:r><w:r><w:t xml:space="preserve">in </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Spring</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> 2018</w:t></w:r><w:r><w:t xml:space="preserve">, and </w:t></w:r><w:r><w:t xml:space="preserve">C1C </w:t></w:r><w:r><w:t xml:space="preserve">Spencer Adolph joined over the summer to help develop the game further. Now as a </w:t></w:r><w:r><w:t>Software Engineering</w:t></w:r><w:r><w:t xml:space="preserve"> capstone</w:t></w:r><w:r><w:t xml:space="preserve"> project</w:t></w:r><w:r><w:t xml:space="preserve">, four </w:t></w:r><w:r><w:t>First Class Computer Science majors</w:t></w:r><w:r><w:t xml:space="preserve"> are developing this game </w:t></w:r><w:r><w:t xml:space="preserve">using </w:t></w:r><w:r><w:t>HTML, JavaScript, and CSS with PHP connection</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> to a MySQL server. The</w:t></w:r><w:r><w:t>ir</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">web-based software </w:t></w:r><w:r><w:t xml:space="preserve">can now handle </w:t></w:r><w:r><w:t xml:space="preserve">multiple </w:t></w:r><w:r><w:t xml:space="preserve">games per section, and automatically </w:t></w:r><w:r><w:t xml:space="preserve">tracks the play of each game using their </w:t></w:r><w:r><w:t xml:space="preserve">database. </w:t></w:r><w:r><w:t xml:space="preserve">MSS 251 cadets </w:t></w:r><w:r><w:t>will be able to drag pieces around the board and initiate combat with a few clicks. The game will have the same phases and turns as the board game.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve"> We plan to have a fully playable version for </w:t></w:r><w:r><w:t>a L</w:t></w:r><w:r><w:t>esson 30 deployment</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>in the MSS</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>251 course</w:t></w:r><w:r><w:t xml:space="preserve"> in the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Fall</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> of 2018</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Replace the whole body paragraph (paragraph 2 - the narrative text) in
# one shot via InsertXML so the tracked insertion (w:ins / rPrChange /
# highlight) and proofErr / bookmark markers land exactly as authored.
$bodyPara = $d.Paragraphs(2).Range
$bodyPara.InsertXML($targetXml)

# The trailing paragraph loses its (now pointless) _GoBack bookmark.
$goBack = $d.Bookmarks("_GoBack")
if ($goBack -ne $null) {
    $goBack.Delete()
}

Write-Output ("Paragraphs=" + $d.Paragraphs.Count)
Write-Output $d.Content.Text
